# Add the automatic repeater scanner scan results: a new PriceChange/UpDown
# pair for the existing row 6, plus a brand new row 7 of scan data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 6 gains a PriceChange (X) / UpDown (Y) verdict.
$ws.Range("X6").Value = 0.29000100000000373
$ws.Range("Y6").Value = "Up"

# New row 7: the next day's scan.
$ws.Range("A7").Value = 42648.890601851854
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "Neutral"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = "Random"
$ws.Range("Q7").Value = 38.916275631518758
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = -0.0078
$ws.Range("S7").NumberFormat = "0.00%"
$ws.Range("T7").Value = -0.0305
$ws.Range("T7").NumberFormat = "0.00%"
$ws.Range("U7").Value = 14.62
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = -2
